$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (crudeoil)
$ws.Range("C2").Value = 546.448087431694
$ws.Range("D2").Value = 34.73848555815769

# Update row 3 values (hydrogen)
$ws.Range("C3").Value = 585.480093676815

# Remove the now-unused "water" and "methane" rows entirely
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()
